# Auto-generated edits applying the Phantom_Profits.xlsx diff
# (sheet tab names ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR correspond to
# the 8 worksheets referenced by the diff as "Sheets/Phantom_Profits.xlsx")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2718.6
$ws.Range("J17").Value = 2910.6667
$ws.Range("L17").Value = 8732.000100000001
$ws.Range("N17").Value = -9068.000100000001

# Row 80
$ws.Range("H80").Value = 1319
$ws.Range("J80").Value = 1556.4166
$ws.Range("L80").Value = 4669.2498
$ws.Range("N80").Value = -6665.2498

# Row 83
$ws.Range("H83").Value = 1319
$ws.Range("J83").Value = 1556.4166
$ws.Range("L83").Value = 14007.7494
$ws.Range("N83").Value = -23991.7494

# Row 97
$ws.Range("H97").Value = 2714.25
$ws.Range("J97").Value = 2714.25
$ws.Range("L97").Value = 8142.75
$ws.Range("N97").Value = -9134.75

# Row 137
$ws.Range("H137").Value = 3001.9412
$ws.Range("I137").Value = 1253.4
$ws.Range("J137").Value = 5499.857
$ws.Range("K137").Value = 3760.2
$ws.Range("L137").Value = 16499.571
$ws.Range("M137").Value = -1210.2
$ws.Range("N137").Value = -21599.571

# Row 138
$ws.Range("H138").Value = 1954.36
$ws.Range("I138").Value = 1779.2727
$ws.Range("J138").Value = 3238.3333
$ws.Range("K138").Value = 5337.8181
$ws.Range("L138").Value = 9714.999899999999
$ws.Range("M138").Value = -197.8181000000004
$ws.Range("N138").Value = -19994.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3263.3125
$ws.Range("I32").Value = 2594.258
$ws.Range("K32").Value = 2594.258
$ws.Range("M32").Value = -2307.258

# Row 37
$ws.Range("H37").Value = 11220.4
$ws.Range("I37").Value = 11220.4
$ws.Range("K37").Value = 11220.4
$ws.Range("M37").Value = -10947.4

# Row 61
$ws.Range("H61").Value = 6189.4194
$ws.Range("I61").Value = 5495.1924
$ws.Range("J61").Value = 9799.4
$ws.Range("K61").Value = 5495.1924
$ws.Range("L61").Value = 9799.4
$ws.Range("M61").Value = -5283.1924
$ws.Range("N61").Value = -10223.4

# Row 110
$ws.Range("H110").Value = 5331.1665
$ws.Range("I110").Value = 5737.7144
$ws.Range("J110").Value = 3908.25
$ws.Range("K110").Value = 5737.7144
$ws.Range("L110").Value = 3908.25
$ws.Range("M110").Value = -3692.7144
$ws.Range("N110").Value = -7998.25

# Row 122
$ws.Range("H122").Value = 1856.8125
$ws.Range("I122").Value = 1737
$ws.Range("J122").Value = 2695.5
$ws.Range("K122").Value = 5211
$ws.Range("L122").Value = 8086.5
$ws.Range("M122").Value = -2761
$ws.Range("N122").Value = -12986.5

# Row 136
$ws.Range("H136").Value = 6189.4194
$ws.Range("I136").Value = 5495.1924
$ws.Range("J136").Value = 9799.4
$ws.Range("K136").Value = 16485.5772
$ws.Range("L136").Value = 29398.2
$ws.Range("M136").Value = -13935.5772
$ws.Range("N136").Value = -34498.2

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 69999
$ws.Range("J35").Value = 69999
$ws.Range("L35").Value = 69999
$ws.Range("N35").Value = -70619

# Row 105
$ws.Range("H105").Value = 5041.8
$ws.Range("I105").Value = 5041.8
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5041.8
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3294.8
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 6474.5
$ws.Range("J58").Value = 7950
$ws.Range("L58").Value = 7950
$ws.Range("N58").Value = -8356

# Row 99
$ws.Range("H99").Value = 3111.5454
$ws.Range("I99").Value = 2778.625
$ws.Range("K99").Value = 2778.625
$ws.Range("M99").Value = -1280.625

# Row 122
$ws.Range("H122").Value = 9599.5625
$ws.Range("I122").Value = 9772.866
$ws.Range("K122").Value = 29318.598
$ws.Range("M122").Value = -26868.598

# Row 126
$ws.Range("H126").Value = 3111.5454
$ws.Range("I126").Value = 2778.625
$ws.Range("K126").Value = 8335.875
$ws.Range("M126").Value = -5865.875

# Row 132
$ws.Range("H132").Value = 1734.5
$ws.Range("J132").Value = 2514
$ws.Range("L132").Value = 7542
$ws.Range("N132").Value = -12602

# Row 136
$ws.Range("H136").Value = 6474.5
$ws.Range("J136").Value = 7950
$ws.Range("L136").Value = 23850
$ws.Range("N136").Value = -28950

$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 1100
$ws.Range("I74").Value = 1100
$ws.Range("K74").Value = 3300
$ws.Range("M74").Value = -2239

# Row 77
$ws.Range("H77").Value = 1100
$ws.Range("I77").Value = 1100
$ws.Range("K77").Value = 9900
$ws.Range("M77").Value = -4596

# Row 107
$ws.Range("H107").Value = 123
$ws.Range("I107").Value = 110
$ws.Range("J107").Value = 129.5
$ws.Range("K107").Value = 330
$ws.Range("L107").Value = 388.5
$ws.Range("M107").Value = 1590
$ws.Range("N107").Value = -4228.5

# Row 131
$ws.Range("H131").Value = 2615
$ws.Range("I131").Value = 2776.5
$ws.Range("K131").Value = 8329.5
$ws.Range("M131").Value = -3289.5

# Row 132
$ws.Range("H132").Value = 1483.5714
$ws.Range("I132").Value = 1328.6666
$ws.Range("K132").Value = 11957.9994
$ws.Range("M132").Value = -9427.999400000001

# Row 138
$ws.Range("H138").Value = 9608.429
$ws.Range("I138").Value = 9608.429
$ws.Range("K138").Value = 28825.287
$ws.Range("M138").Value = -23685.287

# Row 140
$ws.Range("H140").Value = 837712.2
$ws.Range("I140").Value = 913549.6
$ws.Range("K140").Value = 2740648.8
$ws.Range("M140").Value = -2735468.8

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 6490.3335
$ws.Range("I107").Value = 3314.6667
$ws.Range("J107").Value = 9666
$ws.Range("K107").Value = 3314.6667
$ws.Range("L107").Value = 9666
$ws.Range("M107").Value = -1394.6667
$ws.Range("N107").Value = -13506

# Row 122
$ws.Range("H122").Value = 3678.36
$ws.Range("I122").Value = 3627.75
$ws.Range("J122").Value = 3880.8
$ws.Range("K122").Value = 10883.25
$ws.Range("L122").Value = 11642.4
$ws.Range("M122").Value = -8433.25
$ws.Range("N122").Value = -16542.4

# Row 132
$ws.Range("H132").Value = 2513.25
$ws.Range("I132").Value = 2351
$ws.Range("K132").Value = 7053
$ws.Range("M132").Value = -4523

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 999.8
$ws.Range("I16").Value = 999.8
$ws.Range("K16").Value = 999.8
$ws.Range("M16").Value = -829.8

# Row 61
$ws.Range("H61").Value = 3749
$ws.Range("I61").Value = 4498
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 4498
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -4296
$ws.Range("N61").Value = -3404

# Row 68
$ws.Range("H68").Value = 2725
$ws.Range("I68").Value = 2649.5
$ws.Range("K68").Value = 2649.5
$ws.Range("M68").Value = -1900.5

# Row 71
$ws.Range("H71").Value = 2725
$ws.Range("I71").Value = 2649.5
$ws.Range("K71").Value = 13247.5
$ws.Range("M71").Value = -9503.5

# Row 113
$ws.Range("H113").Value = 3749
$ws.Range("I113").Value = 4498
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 4498
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2328
$ws.Range("N113").Value = -7340

# Row 132
$ws.Range("H132").Value = 2486.3
$ws.Range("J132").Value = 2632.6667
$ws.Range("L132").Value = 7898.000100000001
$ws.Range("N132").Value = -12958.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 45749.152
$ws.Range("I45").Value = 38768
$ws.Range("J45").Value = 51733
$ws.Range("K45").Value = 38768
$ws.Range("L45").Value = 51733
$ws.Range("M45").Value = -38277
$ws.Range("N45").Value = -52715

# Row 51
$ws.Range("H51").Value = 184995.78
$ws.Range("I51").Value = 219866
$ws.Range("K51").Value = 219866
$ws.Range("M51").Value = -219356

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 132
$ws.Range("H132").Value = 2517.6177
$ws.Range("J132").Value = 5099.5
$ws.Range("L132").Value = 15298.5
$ws.Range("N132").Value = -20358.5
